# Suite_Web.xlsx - "Suite" sheet update
# Reorders the test-suite rows (column A) and flips most Runmode flags
# (column C) from "N" to "Y", while keeping Web_SETTINGS and
# Web_COMPOSESCREENVALIDATION as "N". Also normalizes the empty
# "placeholder" centered cells in column B to match the new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: reordered test-suite names ---
$ws.Range("A2").Value = "Web_SIGNIN"
$ws.Range("A3").Value = "Web_SEARCH"
$ws.Range("A4").Value = "Web_QUICK_MESSAGES"
$ws.Range("A5").Value = "Web_CONVERSATIONS"
$ws.Range("A6").Value = "Web_PATIENT"
$ws.Range("A7").Value = "Web_URGENT_MESSAGES"
$ws.Range("A8").Value = "Web_GROUPS"
$ws.Range("A9").Value = "Web_USERPROFILE"
$ws.Range("A10").Value = "Web_CONTACTS"
$ws.Range("A11").Value = "Web_SETTINGS"
$ws.Range("A12").Value = "Web_MESSAGING"
$ws.Range("A13").Value = "Web_FILE_UPLOAD_NEW_USER"
$ws.Range("A14").Value = "Web_FILE_UPLOAD_EXISTING_USER"
$ws.Range("A15").Value = "Web_COMPOSESCREENVALIDATION"

# --- Column C: Runmode flags ---
$ws.Range("C2").Value = "Y"
$ws.Range("C3").Value = "Y"
$ws.Range("C4").Value = "Y"
$ws.Range("C5").Value = "Y"
$ws.Range("C6").Value = "Y"
$ws.Range("C7").Value = "Y"
$ws.Range("C8").Value = "Y"
$ws.Range("C9").Value = "Y"
$ws.Range("C10").Value = "Y"
$ws.Range("C11").Value = "N"
$ws.Range("C12").Value = "Y"
$ws.Range("C13").Value = "Y"
$ws.Range("C14").Value = "Y"
$ws.Range("C15").Value = "N"

# --- Column B: empty centered placeholder cells follow the rows that
# now hold the suites which previously occupied B2,B3,B4,B7,B13,B14 ---
$ws.Range("B5").Clear()
$ws.Range("B6").Clear()
$ws.Range("B8").Clear()
$ws.Range("B12").HorizontalAlignment = -4108
$ws.Range("B13").HorizontalAlignment = -4108
$ws.Range("B14").HorizontalAlignment = -4108

# --- Selection moved to C8 (last edited cell) ---
$ws.Range("C8").Select() | Out-Null

# --- Page setup (portrait) ---
$ws.PageSetup.Orientation = 1

